# Add a totals row ("planning and report") to each of the three LogTime
# tables, then restore the selection/scroll state each sheet ended up with.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # ANLT   -> Table24
$ws2 = $wb.Worksheets.Item(2)   # QUANGD -> Table2
$ws3 = $wb.Worksheets.Item(3)   # ANHDT  -> Table22

# ---------------------------------------------------------------------
# Sheet 1 / ANLT / Table24  (A1:M8 -> A1:M9)
# ---------------------------------------------------------------------
$lo1 = $ws1.ListObjects.Item(1)
$lo1.ShowTotals = $true

$ws1.Range("B9").Formula = "=SUM(Table24[10/03/2014])"
$ws1.Range("C9").Formula = "=SUM(Table24[11/03/2014])"
$ws1.Range("D9").Formula = "=SUM(Table24[12/03/2014])"
$ws1.Range("E9").Formula = "=SUM(Table24[13/03/2014])"
$ws1.Range("F9").Formula = "=SUM(Table24[14/03/2014])"
$ws1.Range("G9").Formula = "=SUM(Table24[15/03/2014])"
$ws1.Range("H9").Formula = "=SUM(Table24[17/03/2014])"
$ws1.Range("I9").Formula = "=SUM(Table24[18/03/2014])"
$ws1.Range("J9").Formula = "=SUM(Table24[19/03/2014])"
$ws1.Range("K9").Formula = "=SUM(Table24[20/03/2014])"
$ws1.Range("L9").Formula = "=SUM(Table24[21/03/2014])"
$ws1.Range("M9").Formula = "=SUM(Table24[22/03/2014])"

# ---------------------------------------------------------------------
# Sheet 2 / QUANGD / Table2  (A1:N9 -> A1:N10)
# ---------------------------------------------------------------------
$lo2 = $ws2.ListObjects.Item(1)
$lo2.ShowTotals = $true

$ws2.Range("B10").Formula = "=SUM(Table2[10/03/2014])"
$ws2.Range("C10").Formula = "=SUM(Table2[11/03/2014])"
$ws2.Range("D10").Formula = "=SUM(Table2[12/03/2014])"
$ws2.Range("E10").Formula = "=SUM(Table2[13/03/2014])"
$ws2.Range("F10").Formula = "=SUM(Table2[14/03/2014])"
$ws2.Range("G10").Formula = "=SUM(Table2[15/03/2014])"
$ws2.Range("H10").Formula = "=SUM(Table2[16/03/2015])"
$ws2.Range("I10").Formula = "=SUM(Table2[17/03/2014])"
$ws2.Range("J10").Formula = "=SUM(Table2[18/03/2014])"
$ws2.Range("K10").Formula = "=SUM(Table2[19/03/2014])"
$ws2.Range("L10").Formula = "=SUM(Table2[20/03/2014])"
$ws2.Range("M10").Formula = "=SUM(Table2[21/03/2014])"
$ws2.Range("N10").Formula = "=SUM(Table2[22/03/2014])"

# ---------------------------------------------------------------------
# Sheet 3 / ANHDT / Table22  (A1:N9 -> A1:N10)
# ---------------------------------------------------------------------
$lo3 = $ws3.ListObjects.Item(1)
$lo3.ShowTotals = $true

$ws3.Range("B10").Formula = "=SUM(Table22[10/03/2014])"
$ws3.Range("C10").Formula = "=SUM(Table22[11/03/2014])"
$ws3.Range("D10").Formula = "=SUM(Table22[12/03/2014])"
$ws3.Range("E10").Formula = "=SUM(Table22[13/03/2014])"
$ws3.Range("F10").Formula = "=SUM(Table22[14/03/2014])"
$ws3.Range("G10").Formula = "=SUM(Table22[15/03/2014])"
$ws3.Range("H10").Formula = "=SUM(Table22[16/03/2015])"
$ws3.Range("I10").Formula = "=SUM(Table22[17/03/2014])"
$ws3.Range("J10").Formula = "=SUM(Table22[18/03/2014])"
$ws3.Range("K10").Formula = "=SUM(Table22[19/03/2014])"
$ws3.Range("L10").Formula = "=SUM(Table22[20/03/2014])"
$ws3.Range("M10").Formula = "=SUM(Table22[21/03/2014])"
$ws3.Range("N10").Formula = "=SUM(Table22[22/03/2014])"

# ---------------------------------------------------------------------
# Selections / active sheet, matching the saved view state in the diff.
# Activating sheets in this order leaves ANHDT (sheet3) as the tab shown,
# same as the original workbook (activeTab = 2).
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("K9").Select()

$ws2.Activate()
$ws2.Range("I10:K10").Select()

$ws3.Activate()
$ws3.Range("B10:N10").Select()
